# Lattice-multiplication exercise table: replace the 15 problems (5 rows x 3 cols)
# with a new set of operands/partial-product digits, cell by cell, while preserving
# the original run formatting (sz=32) and line layout (5 lines joined by <w:br/>).
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$wordmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Each entry: row, col, "top operands", "line2 (space-padded)", "line3 dashes", "line4", "line5"
$problems = @(
    @(1, 1, "81 x 94", "  9    4", "  ----", "8|    |", "1|    |"),
    @(1, 2, "95 x 74", "  7    4", "  ----", "9|    |", "5|    |"),
    @(1, 3, "53 x 63", "  6    3", "  ----", "5|    |", "3|    |"),
    @(2, 1, "42 x 40", "  4    0", "  ----", "4|    |", "2|    |"),
    @(2, 2, "17 x 62", "  6    2", "  ----", "1|    |", "7|    |"),
    @(2, 3, "81 x 86", "  8    6", "  ----", "8|    |", "1|    |"),
    @(3, 1, "63 x 33", "  3    3", "  ----", "6|    |", "3|    |"),
    @(3, 2, "13 x 19", "  1    9", "  ----", "1|    |", "3|    |"),
    @(3, 3, "92 x 19", "  1    9", "  ----", "9|    |", "2|    |"),
    @(4, 1, "99 x 94", "  9    4", "  ----", "9|    |", "9|    |"),
    @(4, 2, "59 x 29", "  2    9", "  ----", "5|    |", "9|    |"),
    @(4, 3, "11 x 82", "  8    2", "  ----", "1|    |", "1|    |"),
    @(5, 1, "83 x 32", "  3    2", "  ----", "8|    |", "3|    |"),
    @(5, 2, "76 x 48", "  4    8", "  ----", "7|    |", "6|    |"),
    @(5, 3, "70 x 74", "  7    4", "  ----", "7|    |", "0|    |")
)

foreach ($p in $problems) {
    $row = $p[0]
    $col = $p[1]
    $top = $p[2]
    $line2 = $p[3]
    $line3 = $p[4]
    $line4 = $p[5]
    $line5 = $p[6]
    $cell = $tbl.Cell($row, $col)
    $xml = "<w:p $wordmlNs><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr>" +
           "<w:t>$top</w:t><w:br/>" +
           "<w:t xml:space=`"preserve`">$line2</w:t><w:br/>" +
           "<w:t xml:space=`"preserve`">$line3</w:t><w:br/>" +
           "<w:t>$line4</w:t><w:br/>" +
           "<w:t>$line5</w:t></w:r></w:p>"
    [void]$cell.Range.InsertXML($xml)
}
